$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B (Buying Opportunity) and column C (support Zone) values
# for rows 2-22 with the refreshed ticker list.
$ws.Range("B2").Value = "NSE:APEX"
$ws.Range("C2").Value = "NSE:ATAM"

$ws.Range("B3").Value = "NSE:ARIHANTCAP"
$ws.Range("C3").Value = "NSE:ATGL"

$ws.Range("C4").Value = "NSE:BLS"

$ws.Range("B5").Value = "NSE:BIOCON"
$ws.Range("C5").Value = "NSE:BSHSL"

$ws.Range("B6").Value = "NSE:CHEMFAB"
$ws.Range("C6").Value = "NSE:CAMS"

$ws.Range("B7").Value = "NSE:COLPAL"
$ws.Range("C7").Value = "NSE:CREST"

$ws.Range("B8").Value = "NSE:DONEAR"
$ws.Range("C8").Value = "NSE:CROWN"

$ws.Range("B9").Value = "NSE:FINEORG"
$ws.Range("C9").Value = "NSE:DCM"

$ws.Range("B10").Value = "NSE:GOKEX"
$ws.Range("C10").Value = "NSE:ECLERX"

$ws.Range("B11").Value = "NSE:GREENPANEL"
$ws.Range("C11").Value = "NSE:GENUSPOWER"

$ws.Range("B12").Value = "NSE:GUJGASLTD"
$ws.Range("C12").Value = "NSE:GLAXO"

$ws.Range("B13").Value = "NSE:HDFCLIQUID"
$ws.Range("C13").Value = "NSE:GODFRYPHLP"

$ws.Range("B14").Value = "NSE:IGARASHI"
$ws.Range("C14").Value = "NSE:GREENLAM"

$ws.Range("B15").Value = "NSE:IRMENERGY"
$ws.Range("C15").Value = "NSE:JBCHEPHARM"

$ws.Range("B16").Value = "NSE:MANGALAM"
$ws.Range("C16").Value = "NSE:JINDRILL"

$ws.Range("B17").Value = "NSE:NITCO"
$ws.Range("C17").Value = "NSE:MEDPLUS"

$ws.Range("B18").Value = "NSE:POLYCAB"
$ws.Range("C18").Value = "NSE:MOLDTECH"

$ws.Range("B19").Value = "NSE:PRIVISCL"
$ws.Range("C19").Value = "NSE:MSTCLTD"

$ws.Range("B20").Value = "NSE:RAMCOCEM"
$ws.Range("C20").Value = "NSE:REMSONSIND"

$ws.Range("B21").Value = "NSE:RUSTOMJEE"
$ws.Range("C21").Value = "NSE:RGL"

$ws.Range("B22").Value = "NSE:SAIL"

# Remove the now-obsolete rows 23-36 (tail of the old, longer list).
$ws.Rows("23:36").Delete()
